# The document's headers/footers each carry a small logo picture. Three
# inline pictures need their drawing-object "name" updated:
#   - the two Pearson Edexcel logos in the footers: image1.png -> image2.png
#   - the BTec logo in the (first-page) header:    image2.jpg -> image1.jpg
#
# InlineShape has no settable .Name in the Word object model, so the shape
# is temporarily converted to a floating Shape (which does expose .Name,
# backed by <wp:docPr name="...">), renamed, then converted back to an
# inline shape so the layout/anchoring is restored exactly as it was.

$d = $word.ActiveDocument
$sec = $d.Sections.First

function Rename-InlinePicture($inlineShape, $newName) {
    $floating = $inlineShape.ConvertToShape()
    $floating.Name = $newName
    $floating.ConvertToInlineShape() | Out-Null
}

# --- Footers: both Pearson Edexcel logos, image1.png -> image2.png ---
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $footer = $sec.Footers.Item($i)
    if ($footer.Exists) {
        $shapes = $footer.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                Rename-InlinePicture $shp "image2.png"
            }
        }
    }
}

# --- Headers: BTec logo, image2.jpg -> image1.jpg ---
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $header = $sec.Headers.Item($i)
    if ($header.Exists) {
        $shapes = $header.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                Rename-InlinePicture $shp "image1.jpg"
            }
        }
    }
}

Write-Output "done"
